# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 71, pushing the existing
# rows 71-97 down to 72-98 (dimension grows from A1:R97 to A1:R98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 71 (shifts rows 71:97 -> 72:98)
$ws.Rows.Item(71).EntireRow.Insert()

# Populate the newly inserted row 71 with the new record
$ws.Cells.Item(71, 1).Value  = 2
$ws.Cells.Item(71, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(71, 3).Value  = "Coquimbo"
$ws.Cells.Item(71, 4).Value  = 44524
$ws.Cells.Item(71, 5).Value  = 4
$ws.Cells.Item(71, 6).Value  = 100112024
$ws.Cells.Item(71, 7).Value  = "Choclo"
$ws.Cells.Item(71, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(71, 9).Value  = "Primera"
$ws.Cells.Item(71, 10).Value = 600
$ws.Cells.Item(71, 11).Value = 19000
$ws.Cells.Item(71, 12).Value = 20000
$ws.Cells.Item(71, 13).Value = 19500
$ws.Cells.Item(71, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(71, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(71, 16).Value = 279
$ws.Cells.Item(71, 17).Value = 70
$ws.Cells.Item(71, 18).Value = "Hortaliza"
